# Generate Report for Handoff
# This script updates the localization-status report: the 69f9321c and
# d60bfa6f rows swap places (69f9321c is now "Ready for handoff" and
# moves above d60bfa6f), and a fresh report timestamp is stamped on the
# rows that were (re)generated in this run (3bd8f39e, and the swapped
# 69f9321c/d60bfa6f rows).

$wb = $excel.ActiveWorkbook

function Set-LinkText($ws, $addr, $text) {
    $ws.Range($addr).Value2 = $text
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 6 (3bd8f39e.md) picked up a new handoff timestamp.
$ws.Range("D6").Value2 = "2016-03-23 02:25:10"

# Row 9 used to be d60bfa6f; it is now 69f9321c, ready for handoff.
Set-LinkText $ws '$A$9' "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.md"
$ws.Range("B9").Value2 = "Ready for handoff"
$ws.Range("C9").Value2 = "Ready for handoff"
$ws.Range("D9").Value2 = "2016-03-23 02:25:10"

# Row 10 used to be 69f9321c; it is now d60bfa6f.
Set-LinkText $ws '$A$10' "d60bfa6f-3de0-4b90-a802-0f0c0a47688a.md"
$ws.Range("B10").Value2 = "Ready for handoff"
$ws.Range("C10").Value2 = "Ready for handoff"
$ws.Range("D10").Value2 = "2016-03-23 02:25:10"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 6 (3bd8f39e.md) handoff datetime refreshed.
$ws.Range("E6").Value2 = "2016-03-23 02:25:06"

# Row 9 used to be d60bfa6f; it is now 69f9321c, ready for handoff.
Set-LinkText $ws '$A$9' "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.md"
$ws.Range("C9").Value2 = "Ready for handoff"
Set-LinkText $ws '$D$9' "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.229886f1f51574ec483b3bf868b13ff072de2820.zh-cn.xlf"
$ws.Range("E9").Value2 = "2016-03-23 02:25:06"

# Row 10 used to be 69f9321c; it is now d60bfa6f.
Set-LinkText $ws '$A$10' "d60bfa6f-3de0-4b90-a802-0f0c0a47688a.md"
$ws.Range("C10").Value2 = "Ready for handoff"
Set-LinkText $ws '$D$10' "d60bfa6f-3de0-4b90-a802-0f0c0a47688a.0d1ccd76fe5a8a01c743abfb813e3e53d27fc50a.zh-cn.xlf"
$ws.Range("E10").Value2 = "2016-03-23 02:25:06"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 6 (3bd8f39e.md) handoff datetime refreshed.
$ws.Range("E6").Value2 = "2016-03-23 02:25:10"

# Row 9 used to be d60bfa6f; it is now 69f9321c, ready for handoff.
Set-LinkText $ws '$A$9' "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.md"
$ws.Range("C9").Value2 = "Ready for handoff"
Set-LinkText $ws '$D$9' "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.229886f1f51574ec483b3bf868b13ff072de2820.de-de.xlf"
$ws.Range("E9").Value2 = "2016-03-23 02:25:10"

# Row 10 used to be 69f9321c; it is now d60bfa6f.
Set-LinkText $ws '$A$10' "d60bfa6f-3de0-4b90-a802-0f0c0a47688a.md"
$ws.Range("C10").Value2 = "Ready for handoff"
Set-LinkText $ws '$D$10' "d60bfa6f-3de0-4b90-a802-0f0c0a47688a.0d1ccd76fe5a8a01c743abfb813e3e53d27fc50a.de-de.xlf"
$ws.Range("E10").Value2 = "2016-03-23 02:25:10"
